$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix columns A and B: A should hold month (1-12), B should hold year (2020) ---
# Previously A had the year and B had the month; swap them for each data row.
for ($row = 2; $row -le 13; $row++) {
    $month = $row - 1
    $ws.Cells.Item($row, 1).Value = $month
    $ws.Cells.Item($row, 2).Value = 2020
}

# --- Add new headers in H1:P1 ---
$ws.Range("H1").Value = "grade_total"
$ws.Range("I1").Value = "grade_distance"
$ws.Range("J1").Value = "grade_visitation"
$ws.Range("K1").Value = "grade_encounters"
$ws.Range("L1").Value = "NEVER"
$ws.Range("M1").Value = "RARELY"
$ws.Range("N1").Value = "SOMETIMES"
$ws.Range("O1").Value = "FREQUENTLY"
$ws.Range("P1").Value = "ALWAYS"

# --- Fill new data columns H:P for every data row with the same values ---
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 1.066
    $ws.Cells.Item($row, 13).Value = 1.088
    $ws.Cells.Item($row, 14).Value = 1.122
    $ws.Cells.Item($row, 15).Value = 1.227
    $ws.Cells.Item($row, 16).Value = 1.496
}
